# Junction_Flooding_348: round row-5 readings to the workbook's 2-decimal
# "custom accuracy", and drop the extra (row 6) sample row now that the
# 1000-row dataset regeneration no longer needs it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Tighten row 5's high-precision readings to 2 decimals ---------
$ws.Range("C5").Value  = 16.11
$ws.Range("D5").Value  = 1.14
$ws.Range("E5").Value  = 46.8
$ws.Range("F5").Value  = 38.36
$ws.Range("H5").Value  = 62.43
$ws.Range("I5").Value  = 26.05
$ws.Range("J5").Value  = 11.86
$ws.Range("L5").Value  = 18.76
$ws.Range("M5").Value  = 19.77
$ws.Range("N5").Value  = 5.36
$ws.Range("O5").Value  = 16.92
$ws.Range("P5").Value  = 23.91
$ws.Range("R5").Value  = 0.32
$ws.Range("S5").Value  = 0.8
$ws.Range("T5").Value  = 250.46
$ws.Range("U5").Value  = 47.25
$ws.Range("V5").Value  = 15.62
$ws.Range("W5").Value  = 31.72
$ws.Range("X5").Value  = 16.85
$ws.Range("Y5").Value  = 2.22
$ws.Range("Z5").Value  = 31.68
$ws.Range("AA5").Value = 13.74
$ws.Range("AB5").Value = 12.24
$ws.Range("AE5").Value = 0.43
$ws.Range("AF5").Value = 56.88
$ws.Range("AG5").Value = 8.87
$ws.Range("AH5").Value = 19.47

# --- 2. Drop the now-unneeded row 6 sample -----------------------------
$ws.Rows(6).Delete()

# --- 3. Columns whose rounded values now fit in a narrower column ------
$narrowCols = @(3, 10, 12, 13, 15, 22, 24, 27, 28, 34)
foreach ($c in $narrowCols) {
    $ws.Columns.Item($c).ColumnWidth = 6.14
}
